$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, pushing the existing rows 61..78 down to 62..79
# (this mirrors the weekly data-logging pattern: a new day's reading is added
# at the top of the Rabanito history and everything else shifts down).
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new day's reading.
$ws.Cells.Item(61, 1).Value = 10
$ws.Cells.Item(61, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(61, 3).Value = "La Araucanía"
$ws.Cells.Item(61, 4).Value = 44785
$ws.Cells.Item(61, 5).Value = 9
$ws.Cells.Item(61, 6).Value = 300000001
$ws.Cells.Item(61, 7).Value = "Rabanito"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 50
$ws.Cells.Item(61, 11).Value = 7000
$ws.Cells.Item(61, 12).Value = 8000
$ws.Cells.Item(61, 13).Value = 7600
$ws.Cells.Item(61, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(61, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(61, 16).Value = 633
$ws.Cells.Item(61, 17).Value = 12
$ws.Cells.Item(61, 18).Value = "Hortaliza"
